# Move the counts from column C into column B (column B was empty/unused),
# then clear out column C, and update the column width / view settings to
# match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 108

# Move values from column C to column B for every row that has a value in C.
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $val = $cCell.Value2
    if ($val -ne $null) {
        $ws.Cells.Item($r, 2).Value2 = $val
    }
}

# Clear column C entirely (values + any column-level formatting).
$ws.Range("C1:C$lastRow").Clear()

# Remove the custom width that had been set on column B and instead set
# bestFit/custom width on column A to match the new layout.
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(1).ColumnWidth = 21.42578125

# Update the selection/view: select the whole of column A, and make sure the
# top-left visible cell resets to A1 (no more scrolled-down view).
$ws.Range("A1:A1048576").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
